$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.020.54'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '3.540.45'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.92'
$ws.Range('E5').Value = '  -2.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '197.10'
$ws.Range('E6').Value = '  +5.63%  '
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.654'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.06'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('E12').Value = '  -2.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.54'
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('D14').Value = '4.103.20'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '601.56'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.26'
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('D17').Value = '70.146.71'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.81'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').Value = '3.537.58'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.995'
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.90'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('E23').Value = '  +3.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '102.36'
$ws.Range('E24').Value = '  -2.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.61'
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.14'
$ws.Range('E26').Value = '  +3.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.58'
$ws.Range('E28').Value = '  -2.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.69'
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.14'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('B31').Value = 'dogwifhat'
$ws.Range('C31').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.32'
$ws.Range('E31').Value = '  +18.50%  '
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.42'
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').Value = '0.0₃0837'
$ws.Range('E35').Value = '  +7.18%  '
$ws.Range('D36').Value = '3.791.87'
$ws.Range('E36').Value = '  +6.91%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  -4.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.64'
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.394'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.66'
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '493.49'
$ws.Range('E42').Value = '  -7.68%  '
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0457'
$ws.Range('E44').Value = '  -2.51%  '
$ws.Range('E45').Value = '  -3.47%  '
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.29'
$ws.Range('E47').Value = '  -2.97%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('E49').Value = '  -4.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000249'
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '129.85'
$ws.Range('E51').Value = '  -4.01%  '
